$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 21

# Copy formatting (and formatting only matters here) from the row above
$ws.Range("A20:T20").Copy($ws.Range("A21:T21"))

$ws.Cells.Item($row, 1).Value = 25703000
$ws.Cells.Item($row, 2).Value = "Ангарский"
$ws.Cells.Item($row, 3).Value = 2020
$ws.Cells.Item($row, 4).Value = -0.1924603174603175
$ws.Cells.Item($row, 5).Value = 1.582208568470965
$ws.Cells.Item($row, 6).Value = 0.44572036023971628
$ws.Cells.Item($row, 7).Value = 0.71987124313159401
$ws.Cells.Item($row, 8).Value = 0.53457260648550653
$ws.Cells.Item($row, 9).Value = 0.38583765831023947
$ws.Cells.Item($row, 10).Value = 0.41388166961521788
$ws.Cells.Item($row, 11).Value = 0.38361929426999031
$ws.Cells.Item($row, 12).Value = 0.15259001660970239
$ws.Cells.Item($row, 13).Value = 0.48072894527230958
$ws.Cells.Item($row, 14).Value = 0.024690665053992481
$ws.Cells.Item($row, 15).Value = 0.001977663383487622
$ws.Cells.Item($row, 16).Value = 0.059456628855640528
$ws.Cells.Item($row, 17).Value = 0.019044478153582001
$ws.Cells.Item($row, 18).Value = 0.01995171202809486
$ws.Cells.Item($row, 19).Value = 0.45189827483450068
$ws.Cells.Item($row, 20).Value = 0.34752365837430033

$ws.Range("B26").Select()
